$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DanhSach")
$ws.Activate()

# The "Gioi Tinh" (Giới Tính) column moves from K to G. Cutting column K and
# inserting it before column G shifts the old G:J block (Ngay Sinh, Ngay Rua
# Toi, Ngay Ruoc Le, Ngay Them Suc) one column to the right, to H:K, carrying
# each column's width along with its data - exactly what the target file
# shows.
$ws.Columns("K").Cut() | Out-Null
$ws.Columns("G").Insert() | Out-Null

# Row 3 becomes a copy of row 2's data for columns B:M (column A stays
# blank). Use Copy/PasteSpecial (not a literal Value assignment) so that
# text that looks like a date/number (e.g. "06-02-2010", "0912345678")
# round-trips as the same literal text instead of being re-interpreted.
$ws.Range("B2:M2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("Q12").Select() | Out-Null
